# Additional table formats #12
#
# Adds a new custom table style "no bands" (styleId "nobands"), based on
# "Scroll Table Normal" (styleId "ScrollTableNormal"), to the document's
# style gallery. This mirrors the "no bands" table quick-style that shows
# up at the end of word/styles.xml right after "Scroll Note".

$d = $word.ActiveDocument

# wdStyleTypeTable = 3
$nobands = $d.Styles.Add("nobands", 3)

$nobands.NameLocal = "no bands"
$nobands.BaseStyle = "ScrollTableNormal"
$nobands.Priority  = 99
